$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44377
$ws.Range("H2").Value = "Cultivar IV Región"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17600
$ws.Range("N2").Value = "$/bandeja 18 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 978
$ws.Range("Q2").Value = 18

# Row 3
$ws.Range("D3").Value = 44398
$ws.Range("H3").Value = "Cultivar IV Región"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("N3").Value = "$/bandeja 18 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 972
$ws.Range("Q3").Value = 18

# Row 4
$ws.Range("D4").Value = 44398
$ws.Range("H4").Value = "Cultivar IV Región"
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("N4").Value = "$/bandeja 18 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 861
$ws.Range("Q4").Value = 18

# Row 5
$ws.Range("D5").Value = 44363
$ws.Range("H5").Value = "Cultivar IV Región"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 140
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("N5").Value = "$/bandeja 18 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 806
$ws.Range("Q5").Value = 18

# Row 6
$ws.Range("D6").Value = 44435
$ws.Range("H6").Value = "Cultivar IV Región"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("N6").Value = "$/bandeja 18 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 972
$ws.Range("Q6").Value = 18

# Row 7
$ws.Range("D7").Value = 44435
$ws.Range("H7").Value = "Cultivar IV Región"
$ws.Range("I7").Value = "Tercera"
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14500
$ws.Range("N7").Value = "$/bandeja 18 kilos"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 806
$ws.Range("Q7").Value = 18

# Row 8
$ws.Range("D8").Value = 44391
$ws.Range("H8").Value = "Cultivar IV Región"
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("N8").Value = "$/bandeja 18 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 861
$ws.Range("Q8").Value = 18

# Row 9
$ws.Range("D9").Value = 44433
$ws.Range("H9").Value = "Cultivar IV Región"
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("N9").Value = "$/bandeja 18 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 972
$ws.Range("Q9").Value = 18

# Row 10
$ws.Range("D10").Value = 44433
$ws.Range("H10").Value = "Cultivar IV Región"
$ws.Range("I10").Value = "Tercera"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("N10").Value = "$/bandeja 18 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 806
$ws.Range("Q10").Value = 18

# Row 11
$ws.Range("D11").Value = 44526
$ws.Range("H11").Value = "Cultivar XV región"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5500
$ws.Range("M11").Value = 5250
$ws.Range("N11").Value = "$/caja 10 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 525
$ws.Range("Q11").Value = 10

# Row 12
$ws.Range("D12").Value = 44526
$ws.Range("H12").Value = "Cultivar XV región"
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 4500
$ws.Range("M12").Value = 4250
$ws.Range("N12").Value = "$/caja 10 kilos"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 425
$ws.Range("Q12").Value = 10

# Row 13
$ws.Range("D13").Value = 44526
$ws.Range("H13").Value = "Cultivar XV región"
$ws.Range("I13").Value = "Tercera"
$ws.Range("J13").Value = 120
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3500
$ws.Range("M13").Value = 3250
$ws.Range("N13").Value = "$/caja 10 kilos"
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 325
$ws.Range("Q13").Value = 10

# Row 14
$ws.Range("D14").Value = 44211
$ws.Range("H14").Value = "Cultivar XV región"
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 4500
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = 4750
$ws.Range("N14").Value = "$/caja 10 kilos"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 475
$ws.Range("Q14").Value = 10

# Row 15
$ws.Range("D15").Value = 44405
$ws.Range("H15").Value = "Cultivar IV Región"
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 140
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17500
$ws.Range("N15").Value = "$/bandeja 18 kilos"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 972
$ws.Range("Q15").Value = 18

# Row 16
$ws.Range("D16").Value = 44221
$ws.Range("H16").Value = "Cultivar XV región"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = 5500
$ws.Range("N16").Value = "$/caja 10 kilos"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 550
$ws.Range("Q16").Value = 10

# Row 17
$ws.Range("D17").Value = 44554
$ws.Range("H17").Value = "Cultivar XV región"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 5500
$ws.Range("N17").Value = "$/caja 10 kilos"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 550
$ws.Range("Q17").Value = 10

# Row 18
$ws.Range("D18").Value = 44454
$ws.Range("H18").Value = "Cultivar IV Región"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19500
$ws.Range("N18").Value = "$/bandeja 18 kilos"
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("P18").Value = 1083
$ws.Range("Q18").Value = 18

# Row 19
$ws.Range("D19").Value = 44412
$ws.Range("H19").Value = "Cultivar IV Región"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 17500
$ws.Range("N19").Value = "$/bandeja 18 kilos"
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 972
$ws.Range("Q19").Value = 18

# Row 20
$ws.Range("D20").Value = 44533
$ws.Range("H20").Value = "Cultivar XV región"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 6500
$ws.Range("N20").Value = "$/caja 10 kilos"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 650
$ws.Range("Q20").Value = 10

# Row 21
$ws.Range("D21").Value = 44533
$ws.Range("H21").Value = "Cultivar XV región"
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 4000
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = 4500
$ws.Range("N21").Value = "$/caja 10 kilos"
$ws.Range("O21").Value = "Región de Arica y Parinacota"
$ws.Range("P21").Value = 450
$ws.Range("Q21").Value = 10
